# Update generated counts on "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7238
$wsExhibit.Range("F4").Value = 124
$wsExhibit.Range("G4").Value = 65
$wsExhibit.Range("F5").Value = 177
$wsExhibit.Range("F7").Value = 95
$wsExhibit.Range("F8").Value = 611

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7238
$wsAll.Range("F5").Value = 124
$wsAll.Range("G5").Value = 65
$wsAll.Range("F6").Value = 177
$wsAll.Range("F9").Value = 95
$wsAll.Range("F10").Value = 611
